# fix(publipostage): Correct status name
#
# The "statut_label" column (B) used the color name "bleu" for the
# "no result / no publication yet" status; rename it to "noir".
# The "statut_name" column (C) wording is tightened from
# "... posté(e)" phrasing to "... postés ou publiés" phrasing.
#
# These strings repeat across many rows, so use a whole-cell
# Find & Replace over the used range instead of touching each cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$whole = [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole

# Longer / more specific strings first so a shorter prefix match
# (e.g. "résultat et / ou publication posté") can't shadow the
# more specific 36/12-month variants. xlWhole match means order is
# actually irrelevant for correctness, but keep it tidy anyway.
$ws.Cells.Replace("bleu", "noir", $whole)

$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", $whole)
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", $whole)
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés", $whole)
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", $whole)
